$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C7").Value = -12.812
$ws.Range("A8").Value = -22.273
$ws.Range("A10").Value = -21.754
$ws.Range("A12").Value = -21.771
$ws.Range("C14").Value = -12.271
$ws.Range("C15").Value = -13.636
$ws.Range("A18").Value = -21.784
$ws.Range("C18").Value = -11.162
$ws.Range("C20").Value = -12.252
$ws.Range("A25").Value = -21.832
$ws.Range("C29").Value = -12.124
$ws.Range("C30").Value = -13.347
$ws.Range("C31").Value = -13.277
$ws.Range("C35").Value = -11.986
$ws.Range("A37").Value = -20.184
$ws.Range("C40").Value = -12.782
$ws.Range("C44").Value = -12.827
$ws.Range("C50").Value = -13.326
$ws.Range("C54").Value = -12.978
$ws.Range("A55").Value = -21.868
$ws.Range("A68").Value = -21.736
$ws.Range("C68").Value = -11.001
$ws.Range("C76").Value = -13.045
$ws.Range("A77").Value = -20.42
$ws.Range("A78").Value = -19.915
$ws.Range("A79").Value = -21.57
$ws.Range("A80").Value = -20.15
$ws.Range("A81").Value = -21.933
$ws.Range("A82").Value = -22.037
$ws.Range("A84").Value = -22.14
$ws.Range("C87").Value = -13.213
$ws.Range("C88").Value = -12.985
$ws.Range("C92").Value = -11.403
$ws.Range("C96").Value = -12.894
$ws.Range("C98").Value = -13.23
$ws.Range("A101").Value = -21.055
$ws.Range("C101").Value = -12.723
$ws.Range("A102").Value = -20.044
$ws.Range("C102").Value = -12.731
